# Add new reserv people
# Appends 6 new "reserve list" rows (name / faculty / telegram handle) to
# the bottom of the existing table on "Лист1", rows 64-69, reusing the same
# column layout as the rows directly above them (A = ФИО, B = факультет,
# C = телеграм), then moves the viewport/selection down to the newly added
# data the way the author's Excel session ended up (active cell E65).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new "reserve" rows -----------------------------------------------
# row, height(pt), Name (A), Faculty (B), Telegram (C)
$rows = @(
    @{ Row = 64; Height = 49.25; Name = "Дулепа Милана Богдановна";        Fac = "Юрфак"; Tg = "@mxllln" },
    @{ Row = 65; Height = 49.25; Name = "Комарова Диана Алексеевна";       Fac = "Юрфак"; Tg = "@wwsdiana" },
    @{ Row = 66; Height = 49.25; Name = "Строкова Дарина Дмитриевна";      Fac = "Юрфак"; Tg = "@darri_d" },
    @{ Row = 67; Height = 37.3;  Name = "Лялюева Таисия Сергеевна";        Fac = "Юрфак"; Tg = "@talymndas" },
    @{ Row = 68; Height = 49.25; Name = "Масленникова Светлана Евгеньевна"; Fac = "Юрфак"; Tg = "@s_vv_et" },
    @{ Row = 69; Height = 37.3;  Name = "Самигуллина Зухра Ринатовна";     Fac = "Юрфак"; Tg = "@wqrexo" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    $cellA = $ws.Cells.Item($rowNum, 1)
    $cellA.Value = $r.Name
    $cellA.WrapText = $true

    $cellB = $ws.Cells.Item($rowNum, 2)
    $cellB.Value = $r.Fac

    $cellC = $ws.Cells.Item($rowNum, 3)
    $cellC.Value = $r.Tg
    $cellC.WrapText = $true

    $ws.Rows.Item($rowNum).RowHeight = $r.Height
}

# --- move selection to match the saved view ----------------------------
$ws.Range("E65").Select()
